$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
$ws1.Range("C2").Value = "$2,142.00"

$ws2 = $wb.Worksheets.Item("Squat Stands")
$ws2.Range("C2").Value = "$1,549.00"

$ws3 = $wb.Worksheets.Item("Leg Extensions")
$ws3.Rows.Item(5).Insert()

$ws3.Range("A5").Value = "Selectorized Seated Leg Curl/Extension"
$ws3.Range("B5").Value = "Stray Dog Strength"
$ws3.Range("C5").Value = "$5,250.00"
$ws3.Range("D5").Value = "USA"
$ws3.Range("E5").Value = "https://shop.straydogstrength.com/cdn/shop/files/2325-RIGHT-RED_eee5d4da-9504-4bb9-b7e3-f98e7e85c231.jpg?v=1743705611&width=823"
$ws3.Range("F5").Value = "https://shop.straydogstrength.com/products/selectorized-seated-leg-curl-extension"

$ws3.Hyperlinks.Add($ws3.Range("E5"), "https://shop.straydogstrength.com/cdn/shop/files/2325-RIGHT-RED_eee5d4da-9504-4bb9-b7e3-f98e7e85c231.jpg?v=1743705611&width=823")
$ws3.Hyperlinks.Add($ws3.Range("F5"), "https://shop.straydogstrength.com/products/selectorized-seated-leg-curl-extension")

$ws3.Range("E5").Style = $ws3.Range("E6").Style
$ws3.Range("F5").Style = $ws3.Range("F6").Style
